# Applies the 28-Jun-2020 23:02 COVID data refresh to the "Pais" sheet of paises.xlsx.
# Rows 4-219 are ranked by "Casos totales" (col B) descending, so updated totals that
# overtake a neighbouring country (Egipto/Suecia, Paraguay/Madagascar/..., Fiyi/Dominica,
# Groenlandia/Islas Malvinas) also swap the two rows country-name cells (col A).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Header timestamp banner (A1)
$ws.Range("A1").Value = "Datos actualizados a 28 de Junio de 2020 a las 23:02"

# Row 4: Estados Unidos
$ws.Range("B4").Value = 2629323
$ws.Range("C4").Value = 32786
$ws.Range("D4").Value = 1088393
$ws.Range("E4").Value = 1412541
$ws.Range("G4").Value = 237
$ws.Range("H4").Value = 128389

# Row 17: Alemania
$ws.Range("B17").Value = 194864
$ws.Range("C17").Value = 175
$ws.Range("E17").Value = 8135
$ws.Range("G17").Value = 3
$ws.Range("H17").Value = 9029

# Row 26: Egipto
$ws.Range("A26").Value = "Egipto"
$ws.Range("B26").Value = 65188
$ws.Range("C26").Value = 1265
$ws.Range("D26").Value = 17539
$ws.Range("E26").Value = 44860
$ws.Range("G26").Value = 81
$ws.Range("H26").Value = 2789

# Row 27: Suecia
$ws.Range("A27").Value = "Suecia"
$ws.Range("B27").Value = 65137
$ws.Range("D27").Value = 0
$ws.Range("E27").Value = 0
$ws.Range("H27").Value = 5280

# Row 53: Israel
$ws.Range("B53").Value = 23755
$ws.Range("C53").Value = 334
$ws.Range("D53").Value = 17074
$ws.Range("E53").Value = 6363

# Row 71: Costa de Marfil
$ws.Range("B71").Value = 9101
$ws.Range("C71").Value = 157
$ws.Range("D71").Value = 3808
$ws.Range("E71").Value = 5227

# Row 74: Uzbekistan
$ws.Range("B74").Value = 7948
$ws.Range("C74").Value = 266
$ws.Range("D74").Value = 5329
$ws.Range("E74").Value = 2597
$ws.Range("G74").Value = 2
$ws.Range("H74").Value = 22

# Row 85: Guinea
$ws.Range("B85").Value = 5342
$ws.Range("C85").Value = 51
$ws.Range("D85").Value = 4282
$ws.Range("E85").Value = 1029
$ws.Range("G85").Value = 1
$ws.Range("H85").Value = 31

# Row 99: Costa Rica
$ws.Range("B99").Value = 3130
$ws.Range("C99").Value = 151
$ws.Range("D99").Value = 1366
$ws.Range("E99").Value = 1749

# Row 108: Paraguay
$ws.Range("A108").Value = "Paraguay"
$ws.Range("B108").Value = 2127
$ws.Range("C108").Value = 185
$ws.Range("D108").Value = 1065
$ws.Range("E108").Value = 1047
$ws.Range("G108").Value = 0
$ws.Range("H108").Value = 15

# Row 109: Madagascar
$ws.Range("A109").Value = "Madagascar"
$ws.Range("B109").Value = 2078
$ws.Range("C109").Value = 73
$ws.Range("D109").Value = 944
$ws.Range("E109").Value = 1116
$ws.Range("G109").Value = 2
$ws.Range("H109").Value = 18

# Row 110: Sri Lanka
$ws.Range("A110").Value = "Sri Lanka"
$ws.Range("B110").Value = 2037
$ws.Range("C110").Value = 4
$ws.Range("D110").Value = 1661
$ws.Range("E110").Value = 365
$ws.Range("H110").Value = 11

# Row 111: Guinea Ecuatorial
$ws.Range("A111").Value = "Guinea Ecuatorial"
$ws.Range("B111").Value = 2001
$ws.Range("C111").Value = 0
$ws.Range("D111").Value = 515
$ws.Range("E111").Value = 1454
$ws.Range("H111").Value = 32

# Row 112: Estado de Palestina
$ws.Range("A112").Value = "Estado de Palestina"
$ws.Range("B112").Value = 1990
$ws.Range("C112").Value = 175
$ws.Range("D112").Value = 447
$ws.Range("E112").Value = 1539
$ws.Range("H112").Value = 4

# Row 113: Sudan del Sur
$ws.Range("A113").Value = "Sudan del Sur"
$ws.Range("B113").Value = 1989
$ws.Range("C113").Value = 47
$ws.Range("D113").Value = 246
$ws.Range("E113").Value = 1707
$ws.Range("H113").Value = 36

# Row 114: Estonia
$ws.Range("A114").Value = "Estonia"
$ws.Range("B114").Value = 1987
$ws.Range("C114").Value = 1
$ws.Range("D114").Value = 1818
$ws.Range("E114").Value = 100
$ws.Range("H114").Value = 69

# Row 121: Zambia
$ws.Range("B121").Value = 1557
$ws.Range("C121").Value = 26
$ws.Range("D121").Value = 1311
$ws.Range("E121").Value = 224

# Row 161: Birmania
$ws.Range("B161").Value = 299
$ws.Range("C161").Value = 6
$ws.Range("D161").Value = 218
$ws.Range("E161").Value = 75

# Row 163: Angola
$ws.Range("B163").Value = 267
$ws.Range("C163").Value = 8
$ws.Range("E163").Value = 175
$ws.Range("G163").Value = 1
$ws.Range("H163").Value = 11

# Row 187: Antigua y Barbuda
$ws.Range("B187").Value = 69
$ws.Range("C187").Value = 4
$ws.Range("E187").Value = 44

# Row 191: San Martin (Parte Francesa)
$ws.Range("B191").Value = 43
$ws.Range("C191").Value = 1
$ws.Range("D191").Value = 37

# Row 205: Fiyi
$ws.Range("A205").Value = "Fiyi"

# Row 206: Dominica
$ws.Range("A206").Value = "Dominica"

# Row 209: Groenlandia
$ws.Range("A209").Value = "Groenlandia"

# Row 210: Islas Malvinas
$ws.Range("A210").Value = "Islas Malvinas"
